$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 234, pushing existing rows 234-239 down to 235-240
$ws.Rows("234:234").Insert()

# Populate the newly inserted row 234 with the new data point (12.01.2021).
# The leading apostrophe forces column A to be stored as plain text instead
# of being auto-converted into a date serial by Excel's input parser; the
# format is cleared right after so the cell keeps the sheet's default style
# (matching every other row in the table).
$ws.Range("A234").Value = "'12.01.2021"
$ws.Range("A234").ClearFormats()

$ws.Range("B234").Value = 75067
$ws.Range("C234").Value = 367111
$ws.Range("D234").Value = 2318
$ws.Range("E234").Value = 56281
$ws.Range("F234").Value = 16468
$ws.Range("G234").Value = 0
